# amend table pm11, add column6
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows with revised E/F/G/H figures ---
$updates = @(
    @{ Row = 112; E = 19; F = 30.65; G = 47; H = 1485 },
    @{ Row = 113; E = 23; F = 37.1;  G = 42; H = 1527 },
    @{ Row = 114; E = 21; F = 33.87; G = 40; H = 1567 },
    @{ Row = 115; E = 20; F = 32.26; G = 44; H = 1611 },
    @{ Row = 116; E = 13; F = 20.97; G = 29; H = 1640 },
    @{ Row = 117; H = 1647 },
    @{ Row = 118; H = 1648 },
    @{ Row = 119; E = 25; F = 39.06; G = 62; H = 1710 },
    @{ Row = 120; H = 1738 },
    @{ Row = 121; E = 23; F = 33.82; G = 39; H = 1777 },
    @{ Row = 122; E = 16; F = 23.53; G = 33; H = 1810 },
    @{ Row = 123; E = 10; F = 14.71; G = 15; H = 1825 },
    @{ Row = 124; E = 8;  F = 11.76; G = 8;  H = 1833 },
    @{ Row = 125; H = 1837 },
    @{ Row = 126; E = 13; F = 18.84; G = 31; H = 1868 },
    @{ Row = 127; E = 17; F = 24.29; G = 24; H = 1892 },
    @{ Row = 128; E = 10; F = 14.29; G = 27; H = 1919 }
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey('E')) { $ws.Cells.Item($r, 5).Value = $u.E }
    if ($u.ContainsKey('F')) { $ws.Cells.Item($r, 6).Value = $u.F }
    if ($u.ContainsKey('G')) { $ws.Cells.Item($r, 7).Value = $u.G }
    if ($u.ContainsKey('H')) { $ws.Cells.Item($r, 8).Value = $u.H }
}

# --- Append two new data rows (129, 130) ---
# Copy the date-format styling from the last existing row (A128) so the
# new date cells (A129, A130) pick up the same number format / style index.
$ws.Range("A128").Copy()
$ws.Range("A129:A130").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A129").Value = 45071
$ws.Range("B129").Value = 21
$ws.Range("C129").Value = 2
$ws.Range("D129").Value = 72
$ws.Range("E129").Value = 18
$ws.Range("F129").Value = 25
$ws.Range("G129").Value = 32
$ws.Range("H129").Value = 1951

$ws.Range("A130").Value = 45072
$ws.Range("B130").Value = 21
$ws.Range("C130").Value = 0
$ws.Range("D130").Value = 72
$ws.Range("E130").Value = 4
$ws.Range("F130").Value = 5.56
$ws.Range("G130").Value = 6
$ws.Range("H130").Value = 1957
